# Fix cell inputs across each property-analysis worksheet.
#
# For every sheet in the workbook:
#  1. Remove the stray column-E "Yearly" duplicate block (E19:E25) that
#     duplicated column D's yearly figures.
#  2. Correct the Repairs/Vacancy/CapEx/Management "(%-mo)" inputs in
#     B22:B25 - they had been entered as absolute dollar amounts instead
#     of the intended percentages.
#  3. Fix the "Annualized Return" formulas in row 34 so the holding
#     period used in the exponent is (Year + 1) instead of Year, and so
#     B34 uses the same profit-based formula as the rest of the row.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # 1) Drop the duplicated column-E "Yearly" cells for rows 19-25.
    $ws.Range("E19:E25").ClearContents()

    # 2) Replace the dollar-value inputs with the intended percentages.
    $ws.Range("B22").Value = 0.05
    $ws.Range("B23").Value = 0.09
    $ws.Range("B24").Value = 0.1
    $ws.Range("B25").Value = 0.1

    # 3) Correct the Annualized Return row formulas.
    $ws.Range("B34").Formula = "=((B33+E6)/E6)^(1/(B27+1))-1"
    $ws.Range("C34").Formula = "=((C33+E6)/E6)^(1/(C27+1))-1"
    $ws.Range("D34").Formula = "=((D33+E6)/E6)^(1/(D27+1))-1"
    $ws.Range("E34").Formula = "=((E33+E6)/E6)^(1/(E27+1))-1"
    $ws.Range("F34").Formula = "=((F33+E6)/E6)^(1/(F27+1))-1"
    $ws.Range("G34").Formula = "=((G33+E6)/E6)^(1/(G27+1))-1"
}
